$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25, shifting existing rows 25-30 down to 26-31
$ws.Rows("25:25").Insert()

# Populate the newly inserted row 25 with the new data record
$ws.Range("A25").Value = 8
$ws.Range("B25").Value = "Terminal La Palmera de La Serena"
$ws.Range("C25").Value = "Coquimbo"
$ws.Range("D25").Value = 45124
$ws.Range("E25").Value = 4
$ws.Range("F25").Value = 100112013
$ws.Range("G25").Value = "Alcachofa"
$ws.Range("H25").Value = "Española"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 400
$ws.Range("K25").Value = 14000
$ws.Range("L25").Value = 15000
$ws.Range("M25").Value = 14500
$ws.Range("N25").Value = "$/caja 30 unidades"
$ws.Range("O25").Value = "Provincia de Limarí"
$ws.Range("P25").Value = 483
$ws.Range("Q25").Value = 30
$ws.Range("R25").Value = "Hortaliza"
